$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The reporting week shifted forward by 7 days: every date in column A
# (each date spans a block of 9 task rows) moves from 2023-09-18..23 to
# 2023-09-25..30.
$dateMap = @{
    "2023-09-18" = "2023-09-25"
    "2023-09-19" = "2023-09-26"
    "2023-09-20" = "2023-09-27"
    "2023-09-21" = "2023-09-28"
    "2023-09-22" = "2023-09-29"
    "2023-09-23" = "2023-09-30"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $current = $cell.Value()
    if ($dateMap.ContainsKey($current)) {
        $cell.Value = $dateMap[$current]
    }
}

# Scroll the view up a page and move the active selection accordingly
# (was topLeftCell A43 / selection B47 -> now topLeftCell A34 / selection B42).
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B42").Select()
